$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$words = @(
    'streichen',
    'stehlen',
    'tollen',
    'saufen',
    'rufen',
    'trennen',
    'filmen',
    'bluten',
    'heben',
    'backen',
    'zünden',
    'feiern',
    'töten',
    'jubeln',
    'knarren',
    'rasen',
    'boxen',
    'bitten',
    'arten',
    'sinken',
    'stecken',
    'flüchten',
    'hauen',
    'schreiten',
    'irren',
    'sichern',
    'äußern',
    'schulden',
    'greifen',
    'quälen',
    'fangen',
    'decken',
    'seufzen',
    'helfen',
    'spielen',
    'trauen',
    'dringen',
    'schlucken',
    'platzen',
    'liegen',
    'bauen',
    'sorgen',
    'kümmern',
    'wenden',
    'freuen',
    'wachsen',
    'erben',
    'machen',
    'fahren',
    'sterben',
    'lockern',
    'stammen',
    'kosten',
    'zielen',
    'drehen',
    'pfeifen',
    'mauern',
    'pflanzen',
    'enden',
    'biegen',
    'ändern',
    'reizen',
    'werden',
    'heilen',
    'ehren',
    'fließen',
    'achten',
    'formen',
    'führen',
    'schwören',
    'räumen',
    'runden',
    'scheitern',
    'kehren',
    'brauchen',
    'altern',
    'grüßen',
    'sprengen',
    'klingen',
    'spüren',
    'flehen',
    'fallen',
    'geben',
    'fällen',
    'bergen',
    'dienen',
    'liefern',
    'schwächen',
    'schrecken',
    'suchen',
    'kichern',
    'ärgern',
    'weichen',
    'tropfen',
    'mögen',
    'klettern',
    'lügen',
    'lesen',
    'schenken',
    'münzen',
    'planen',
    'schmecken',
    'treiben',
    'werfen',
    'hören',
    'graben',
    'fischen',
    'zögern',
    'betteln',
    'schwingen',
    'wüten',
    'sperren',
    'bellen',
    'warnen',
    'heulen',
    'gelten',
    'siegen',
    'folgen',
    'zeigen',
    'scheinen',
    'wundern',
    'jagen',
    'spinnen',
    'wehtun',
    'loben',
    'gründen',
    'wirken',
    'malen'
)
for ($i = 0; $i -lt $words.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $words[$i]
}
